# Generate Report for Handback
# The f66ea9a6-bd9e-40fe-a23d-619da94fec1c file has been handed back and is
# in sync with en-US. Update its status on every sheet and stamp the
# handback datetime for each locale.

$wb = $excel.ActiveWorkbook

$status = "Handed back: in sync with en-US"

# Overview sheet: update both locale status columns for the f66ea9a6... row
$wsOverview = $wb.Worksheets.Item("Overview")
$wsOverview.Range("B3").Value = $status
$wsOverview.Range("C3").Value = $status

# zh-cn sheet: update Status column and Latest Handback DateTime
$wsZhCn = $wb.Worksheets.Item("zh-cn")
$wsZhCn.Range("C3").Value = $status
$wsZhCn.Range("H3").Value = "2016-03-31 07:10:42"

# de-de sheet: update Status column and Latest Handback DateTime
$wsDeDe = $wb.Worksheets.Item("de-de")
$wsDeDe.Range("C3").Value = $status
$wsDeDe.Range("H3").Value = "2016-03-31 07:11:01"
